$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D width (4th column) so the saved OOXML width attribute is 40
$ws.Columns.Item(4).ColumnWidth = 39.1

# Update row 2 values
$ws.Range("A2").Value = "Você"
$ws.Range("B2").Value = "Desconhecido"
$ws.Range("C2").Value = "R$ 200,00"
$ws.Range("D2").Value = "COMERCIO DE POLPAS SOUZA E DIAS LTD..."

# Update row 3 value
$ws.Range("B3").Value = "20:28, 25/03/2025"
